$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.847.09"
$ws.Range("E2").Value = "  -4.32%  "

$ws.Range("D3").Value = "2.229.31"
$ws.Range("E3").Value = "  -5.90%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "491.26"
$ws.Range("E5").Value = "  -2.96%  "

$ws.Range("D6").Value = "127.11"
$ws.Range("E6").Value = "  -2.35%  "

$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("E8").Value = "  -2.94%  "

$ws.Range("D9").Value = "2.266.29"
$ws.Range("E9").Value = "  -4.65%  "

$ws.Range("D10").Value = "0.0930"
$ws.Range("E10").Value = "  -5.78%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").Value = "4.65"
$ws.Range("E13").Value = "  -4.78%  "

$ws.Range("D14").Value = "2.628.99"
$ws.Range("E14").Value = "  -5.76%  "

$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "53.793.72"
$ws.Range("E16").Value = "  -4.19%  "

$ws.Range("E17").Value = "  -3.45%  "

$ws.Range("D18").Value = "2.242.00"
$ws.Range("E18").Value = "  -5.25%  "

$ws.Range("E19").Value = "  -2.39%  "

$ws.Range("D20").Value = "4.04"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").Value = "297.66"
$ws.Range("E21").Value = "  -3.92%  "

$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.40%  "

$ws.Range("D24").Value = "63.53"
$ws.Range("E24").Value = "  -3.26%  "

$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("D26").Value = "0.374"
$ws.Range("E26").Value = "  +0.68%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.147"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.330.15"
$ws.Range("E28").Value = "  -6.13%  "

$ws.Range("D29").Value = "7.12"
$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("D30").Value = "162.72"
$ws.Range("E30").Value = "  -5.97%  "

$ws.Range("D31").Value = "1.60"
$ws.Range("E31").Value = "  -2.94%  "

$ws.Range("D32").Value = "0.0₃0679"
$ws.Range("E32").Value = "  -4.56%  "

$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").Value = "5.81"
$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("D35").Value = "0.991"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").Value = "17.42"
$ws.Range("E37").Value = "  -0.95%  "

$ws.Range("D38").Value = "1.18"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").Value = "0.842"
$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -1.46%  "

$ws.Range("D41").Value = "35.36"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("D42").Value = "0.375"
$ws.Range("E42").Value = "  +0.89%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "4.91"
$ws.Range("E45").Value = "  +3.55%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "125.69"
$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").Value = "0.0890"
$ws.Range("E47").Value = "  -1.01%  "

$ws.Range("D48").Value = "242.57"
$ws.Range("E48").Value = "  +2.20%  "

$ws.Range("D49").Value = "0.543"
$ws.Range("E49").Value = "  -3.43%  "

$ws.Range("E50").Value = "  -0.75%  "

$ws.Range("D51").Value = "0.0203"
$ws.Range("E51").Value = "  -1.83%  "
